# Update the Madigan bike hours (Riders/Average) figures on the Ridership sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# New Riders (column C) and Average (column D) values per row
$ws.Range("C2").Value = 211
$ws.Range("D2").Value = 230.75

$ws.Range("C3").Value = 176
$ws.Range("D3").Value = 209.38

$ws.Range("C4").Value = 204
$ws.Range("D4").Value = 194.92

$ws.Range("C5").Value = 220
$ws.Range("D5").Value = 224.67

$ws.Range("C6").Value = 236
$ws.Range("D6").Value = 237.23

$ws.Range("C7").Value = 111
$ws.Range("D7").Value = 121.14

$ws.Range("C8").Value = 131
$ws.Range("D8").Value = 105.92

$wb.Save()
